$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text for the 93770ace row changes from
#     "Ready for handoff" to "Handback transform failed" (Status + dup column)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

# --- zh-cn sheet: Status column (C) for the 93770ace row, plus a new
#     "Error Detail" message in column P, and widen column P to fit it.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("P3").Value = "Handback file name: f0msrena.zr0 is different with handoff file name: 93770ace-3d57-4dc1-b3ad-0c9cc51301d2.b4b0e60f106c185485ae49c898ddfb3aabc68d46.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: same treatment
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("P3").Value = "Handback file name: f0msrena.zr0 is different with handoff file name: 93770ace-3d57-4dc1-b3ad-0c9cc51301d2.b4b0e60f106c185485ae49c898ddfb3aabc68d46.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
